$d = $word.ActiveDocument

$pairs = @(
    @("49÷8=", "80÷2="),
    @("27÷7=", "14÷2="),
    @("79÷6=", "61÷4="),
    @("51÷3=", "34÷8="),
    @("89÷5=", "28÷6="),
    @("51÷4=", "59÷5="),
    @("10÷7=", "89÷4="),
    @("19÷3=", "40÷9="),
    @("57÷3=", "92÷5="),
    @("55÷2=", "81÷2="),
    @("50÷7=", "70÷3="),
    @("48÷2=", "82÷6="),
    @("21÷2=", "95÷9="),
    @("85÷2=", "80÷9="),
    @("19÷8=", "17÷3="),
    @("16÷6=", "33÷5="),
    @("75÷2=", "59÷8="),
    @("76÷2=", "76÷4="),
    @("64÷4=", "60÷9="),
    @("88÷9=", "24÷7="),
    @("75÷7=", "13÷8="),
    @("55÷6=", "26÷7="),
    @("20÷2=", "42÷7="),
    @("31÷6=", "87÷9="),
    @("59÷7=", "84÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
